$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.905.34'
$ws.Range('E2').Value = '  +3.55%  '
$ws.Range('D3').Value = '1.679.96'
$ws.Range('E3').Value = '  +3.30%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').Value = '219.67'
$ws.Range('E5').Value = '  +2.47%  '
$ws.Range('E6').Value = '  +3.25%  '
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('D8').Value = '29.02'
$ws.Range('E8').Value = '  +2.08%  '
$ws.Range('E9').Value = '  +2.92%  '
$ws.Range('D10').Value = '0.0644'
$ws.Range('E10').Value = '  +5.94%  '
$ws.Range('D11').Value = '0.0904'
$ws.Range('E11').Value = '  +0.37%  '
$ws.Range('D12').Value = '1.921.18'
$ws.Range('E12').Value = '  +3.28%  '
$ws.Range('D13').Value = '1.683.76'
$ws.Range('E13').Value = '  +3.51%  '
$ws.Range('D14').Value = '0.603'
$ws.Range('E14').Value = '  +7.14%  '
$ws.Range('D15').Value = '10.07'
$ws.Range('E15').Value = '  +9.38%  '
$ws.Range('D16').Value = '4.11'
$ws.Range('E16').Value = '  +7.56%  '
$ws.Range('D17').Value = '30.856.05'
$ws.Range('E17').Value = '  +3.30%  '
$ws.Range('D18').Value = '66.00'
$ws.Range('E18').Value = '  +2.33%  '
$ws.Range('D19').Value = '244.46'
$ws.Range('E19').Value = '  +1.98%  '
$ws.Range('D20').Value = '0.0₃0720'
$ws.Range('E20').Value = '  +2.83%  '
$ws.Range('D21').Value = '1.00'
$ws.Range('E21').Value = '  -0.06%  '
$ws.Range('D22').Value = '4.24'
$ws.Range('E22').Value = '  +3.31%  '
$ws.Range('D23').Value = '9.97'
$ws.Range('E23').Value = '  +1.97%  '
$ws.Range('D24').Value = '2.16'
$ws.Range('E24').Value = '  +0.08%  '
$ws.Range('D25').Value = '159.23'
$ws.Range('E25').Value = '  +1.01%  '
$ws.Range('D26').Value = '15.84'
$ws.Range('E26').Value = '  +2.62%  '
$ws.Range('E27').Value = '  +2.73%  '
$ws.Range('D28').Value = '6.68'
$ws.Range('E28').Value = '  +2.21%  '
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.18%  '
$ws.Range('D30').Value = '0.0493'
$ws.Range('E30').Value = '  +1.13%  '
$ws.Range('E31').Value = '  +3.73%  '
$ws.Range('E32').Value = '  +3.26%  '
$ws.Range('D33').Value = '1.523.14'
$ws.Range('E33').Value = '  +7.07%  '
$ws.Range('E34').Value = '  +4.58%  '
$ws.Range('D35').Value = '1.75'
$ws.Range('E35').Value = '  +4.76%  '
$ws.Range('D36').Value = '84.11'
$ws.Range('E36').Value = '  +12.76%  '
$ws.Range('E37').Value = '  +0.72%  '
$ws.Range('D38').Value = '0.605'
$ws.Range('E38').Value = '  +9.00%  '
$ws.Range('E39').Value = '  +5.52%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = '2.66'
$ws.Range('E40').Value = '  -3.40%  '
$ws.Range('B41').Value = 'HuobiToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D41').Value = '2.29'
$ws.Range('E41').Value = '  +0.10%  '
$ws.Range('E42').Value = '  +3.36%  '
$ws.Range('D43').Value = '0.839'
$ws.Range('E43').Value = '  +1.57%  '
$ws.Range('E44').Value = '  +0.37%  '
$ws.Range('E45').Value = '  +2.59%  '
$ws.Range('D46').Value = '1.00'
$ws.Range('E46').Value = '  -0.11%  '
$ws.Range('D47').Value = '5.57'
$ws.Range('E47').Value = '  +4.73%  '
$ws.Range('D48').Value = '50.86'
$ws.Range('E48').Value = '  +5.11%  '
$ws.Range('D49').Value = '1.812.08'
$ws.Range('E49').Value = '  +2.55%  '
$ws.Range('D50').Value = '0.0₆0119'
$ws.Range('E50').Value = '  +7.37%  '
$ws.Range('D51').Value = '92.93'
$ws.Range('E51').Value = '  +2.30%  '
